$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip trailing newline from the player-name header cells.
$ws.Range("A1").Value = "Kim"
$ws.Range("C1").Value = "Emil"
$ws.Range("E1").Value = "Mads"
$ws.Range("G1").Value = "Soren"

# Fix / normalize team-name spelling and casing.
$ws.Range("A3").Value = "RB Leipzig"
$ws.Range("E3").Value = "Eintracht Frankfurt"
$ws.Range("G4").Value = "Sevilla"
$ws.Range("C5").Value = "Real Sociedad"
$ws.Range("G5").Value = "Juventus"
$ws.Range("E6").Value = "Milan"
$ws.Range("C7").Value = "Brøndby"

# Turn the hard-coded 0 score cells into "=0" formulas.
$scoreCells = @("B2","D2","F2","H2","B3","D3","F3","H3","B4","D4","F4","H4","B5","D5","F5","H5","B6","D6","F6","H6","B7","D7","F7","H7")
foreach ($cellRef in $scoreCells) {
    $ws.Range($cellRef).Formula = "=0"
}
